# FUNCTIONALITY: Relocated and wrote out a few test cases.
#
# The external-link source data (Sheet1!H5 = 228, Sheet1!H6 = 189 in the
# linked "Create/_Test_Suite_Statistics_for_folders.xlsx" workbook) was
# bumped by the test-suite run: H5 228 -> 238, H6 189 -> 199.
#
# That feeds Sheet1 of this workbook via the cached external-reference
# formulas in E2 (=[1]Sheet1!$H$6) and F2 (=[1]Sheet1!$H$5), and from there
# into the SUM()/ratio formulas in H5, H6 and H7.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Updated cached results pulled in from the external workbook link.
$ws.Range("F2").Value = 238
$ws.Range("E2").Value = 199

# Recalculate so the dependent SUM/ratio formulas (H5: SUM($F:$F),
# H6: SUM($E:$E), H7: H6/H5) pick up the new inputs and re-cache their
# results, same as Excel does when you update links.
$excel.CalculateFull()
